$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "501.88") must keep
# their original Text storage (the source sheet stores every Price/Volume
# cell as text) so force a Text number format before writing those values.
$textCells = @("D5", "D6", "D10", "D13", "D16", "D19", "D21", "D24", "D26", "D27", "D28", "D29", "D34", "D36", "D38", "D39", "D40", "D42", "D45", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "55.994.17"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.362.97"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "501.88"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "130.77"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("D9").Value = "2.358.02"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").Value = "0.0972"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "4.63"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "2.782.53"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "55.934.69"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Value = "21.36"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "2.376.72"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "9.98"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").Value = "305.80"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "64.94"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").Value = "0.147"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "7.21"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").Value = "172.20"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "0.0₃0714"
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("E35").Value = "  -7.70%  "
$ws.Range("D36").Value = "17.55"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").Value = "3.77"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").Value = "36.00"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "0.793"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "130.78"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0907"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.560"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "242.94"
$ws.Range("E47").Value = "  -5.87%  "
$ws.Range("D48").Value = "0.0478"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("D50").Value = "16.93"
$ws.Range("E50").Value = "  -0.75%  "
